# Generate Report for Handoff
#
# The four files that were still "Ready for handoff" (18414157, 4904a974,
# d05925c5, f2ce4df9) just had their localization xliff handed off again
# (re-generated) for both target languages (zh-cn and de-de). This updates
# the Priority from "low" to "ht", bumps the "Latest Handoff Datetime" to
# the new generation time, fills in the target/handback file columns for
# the two rows that didn't have them yet, and refreshes the rollup
# "Latest HO Xliff Generate Date" on the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$zhCnHandoffTime = "2016-08-16 12:30:47"
$deDeHandoffTime = "2016-08-16 12:30:52"

# --- Overview sheet: rollup "Latest HO Xliff Generate Date" (rows 4-7) ---
$wsOverview.Range("G4").Value = $deDeHandoffTime
$wsOverview.Range("G5").Value = $deDeHandoffTime
$wsOverview.Range("G6").Value = $deDeHandoffTime
$wsOverview.Range("G7").Value = $deDeHandoffTime

# --- zh-cn sheet: Priority + Latest Handoff Datetime for rows 4-7 ---
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("H4").Value = $zhCnHandoffTime

$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("H5").Value = $zhCnHandoffTime

$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("H6").Value = $zhCnHandoffTime

$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("H7").Value = $zhCnHandoffTime

# rows 4 and 5 did not yet have a "Latest Target File" / "Latest Handback
# File" recorded - fill them in now that the handoff went out.
$wsZhCn.Range("I4").Value = "18414157-04ae-426e-839e-f75ea1599d93.md"
$wsZhCn.Range("J4").Value = "18414157-04ae-426e-839e-f75ea1599d93.539e8c7b4869d77fd8fe3aa842a52d54bd1d0b2b.zh-cn.xlf"

$wsZhCn.Range("I5").Value = "4904a974-f06a-423f-8e1b-691bca78932a.md"
$wsZhCn.Range("J5").Value = "4904a974-f06a-423f-8e1b-691bca78932a.90e97ca008ab8e33348807e3cbc45c0c9dc3ef32.zh-cn.xlf"

# --- de-de sheet: Priority + Latest Handoff Datetime for rows 4-7 ---
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("H4").Value = $deDeHandoffTime

$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("H5").Value = $deDeHandoffTime

$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("H6").Value = $deDeHandoffTime

$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("H7").Value = $deDeHandoffTime
